$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 3: end time changed from 11:00 to 10:00 (actual worked time drops to 1.5h)
$ws.Range("E3").Value = 0.41666666666666669

# Row 4: fill in the previously-empty row with a new task entry
$ws.Range("A4").Value = "Review på UC07"
$ws.Range("B4").Value = "Reviewer"
$ws.Range("C4").Value = 43893
$ws.Range("D4").Value = 0.41666666666666669
$ws.Range("E4").Value = 0.45833333333333331
$ws.Range("F4").Value = 0.041666666666666664
$ws.Range("F4").NumberFormat = "h:mm"

# Row 5: fill in the previously-empty row with another new task entry
$ws.Range("A5").Value = "Merge på OC0803"
$ws.Range("B5").Value = "Deployment Manager"
$ws.Range("C5").Value = 43893
$ws.Range("D5").Value = 0.58333333333333337
$ws.Range("E5").Value = 0.625

# Move the active selection to E8 (matches the author's last position when saving)
$ws.Range("E8").Select()
